# Scheduled-runner refresh of market-price-derived profit figures across
# the per-job leve-profit sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Only specific H:N cells (currentAveragePrice.. / LevePrice.. / LeveProfit..)
# on specific rows are refreshed with newly pulled values; everything else
# (names, levels, gil, item ids, etc.) is left untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 61 (ALC)
$ws.Range("H61").Value = 212.5
$ws.Range("I61").Value = 150
$ws.Range("J61").Value = 400
$ws.Range("K61").Value = 450
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -278
$ws.Range("N61").Value = -1544

# Row 112 (ALC)
$ws.Range("H112").Value = 1014.89795
$ws.Range("J112").Value = 1016.08246
$ws.Range("L112").Value = 3048.24738
$ws.Range("N112").Value = -5264.24738

# Row 137 (ALC)
$ws.Range("H137").Value = 1466.3077
$ws.Range("I137").Value = 1067.3334
$ws.Range("J137").Value = 1808.2858
$ws.Range("K137").Value = 3202.0002
$ws.Range("L137").Value = 5424.857400000001
$ws.Range("M137").Value = -652.0001999999999
$ws.Range("N137").Value = -10524.8574

# Row 138 (ALC)
$ws.Range("H138").Value = 5052.8047
$ws.Range("I138").Value = 1798.3334
$ws.Range("J138").Value = 6930.385
$ws.Range("K138").Value = 5395.0002
$ws.Range("L138").Value = 20791.155
$ws.Range("M138").Value = -255.0002000000004
$ws.Range("N138").Value = -31071.155

$ws = $wb.Worksheets.Item("ARM")
# Row 6 (ARM)
$ws.Range("H6").Value = 15000
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

# Row 32 (ARM)
$ws.Range("H32").Value = 5467182.5
$ws.Range("I32").Value = 2407.3076
$ws.Range("J32").Value = 37041440
$ws.Range("K32").Value = 2407.3076
$ws.Range("L32").Value = 37041440
$ws.Range("M32").Value = -2120.3076
$ws.Range("N32").Value = -37042014

# Row 61 (ARM)
$ws.Range("H61").Value = 4631046.5
$ws.Range("I61").Value = 4631046.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4631046.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4630834.5
$ws.Range("N61").ClearContents()

# Row 136 (ARM)
$ws.Range("H136").Value = 4631046.5
$ws.Range("I136").Value = 4631046.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 13893139.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -13890589.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 105 (BSM)
$ws.Range("H105").Value = 200006000
$ws.Range("I105").Value = 7500
$ws.Range("J105").Value = 1000000000
$ws.Range("K105").Value = 7500
$ws.Range("L105").Value = 1000000000
$ws.Range("M105").Value = -5753
$ws.Range("N105").Value = -1000003494

# Row 134 (BSM)
$ws.Range("H134").Value = 3973530
$ws.Range("I134").Value = 1691.9131
$ws.Range("J134").Value = 22243986
$ws.Range("K134").Value = 5075.7393
$ws.Range("L134").Value = 66731958
$ws.Range("M134").Value = -2540.7393
$ws.Range("N134").Value = -66737028

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (CRP)
$ws.Range("H4").Value = 5010000
$ws.Range("J4").Value = 5010000
$ws.Range("L4").Value = 5010000
$ws.Range("N4").Value = -5010224

# Row 58 (CRP)
$ws.Range("H58").Value = 38462204
$ws.Range("I58").Value = 45455176
$ws.Range("J58").Value = 857
$ws.Range("K58").Value = 45455176
$ws.Range("L58").Value = 857
$ws.Range("M58").Value = -45454973
$ws.Range("N58").Value = -1263

# Row 74 (CRP)
$ws.Range("H74").Value = 17325.715
$ws.Range("J74").Value = 16046.667
$ws.Range("L74").Value = 16046.667
$ws.Range("N74").Value = -17794.667

# Row 77 (CRP)
$ws.Range("H77").Value = 17325.715
$ws.Range("J77").Value = 16046.667
$ws.Range("L77").Value = 48140.001
$ws.Range("N77").Value = -56876.001

# Row 122 (CRP)
$ws.Range("H122").Value = 12501376
$ws.Range("I122").Value = 19232324
$ws.Range("J122").Value = 1043.4286
$ws.Range("K122").Value = 57696972
$ws.Range("L122").Value = 3130.2858
$ws.Range("M122").Value = -57694522
$ws.Range("N122").Value = -8030.2858

# Row 134 (CRP)
$ws.Range("H134").Value = 14286685
$ws.Range("I134").Value = 869.9355
$ws.Range("J134").Value = 125001750
$ws.Range("K134").Value = 2609.8065
$ws.Range("L134").Value = 375005250
$ws.Range("M134").Value = -74.80650000000014
$ws.Range("N134").Value = -375010320

# Row 136 (CRP)
$ws.Range("H136").Value = 38462204
$ws.Range("I136").Value = 45455176
$ws.Range("J136").Value = 857
$ws.Range("K136").Value = 136365528
$ws.Range("L136").Value = 2571
$ws.Range("M136").Value = -136362978
$ws.Range("N136").Value = -7671

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (CUL)
$ws.Range("H4").Value = 125100
$ws.Range("I4").Value = 125100
$ws.Range("K4").Value = 375300
$ws.Range("M4").Value = -375188

# Row 99 (CUL)
$ws.Range("H99").Value = 1941.6666
$ws.Range("I99").Value = 1912.5
$ws.Range("K99").Value = 5737.5
$ws.Range("M99").Value = -3491.5

# Row 131 (CUL)
$ws.Range("H131").Value = 911.52
$ws.Range("I131").Value = 607.5
$ws.Range("J131").Value = 924.1875
$ws.Range("K131").Value = 1822.5
$ws.Range("L131").Value = 2772.5625
$ws.Range("M131").Value = 3217.5
$ws.Range("N131").Value = -12852.5625

# Row 134 (CUL)
$ws.Range("H134").Value = 33335380
$ws.Range("I134").Value = 55556468
$ws.Range("J134").Value = 3749.5
$ws.Range("K134").Value = 166669404
$ws.Range("L134").Value = 11248.5
$ws.Range("M134").Value = -166664334
$ws.Range("N134").Value = -21388.5

# Row 139 (CUL)
$ws.Range("H139").Value = 4338.3335
$ws.Range("I139").Value = 1030
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 3090
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = 2050
$ws.Range("N139").Value = -25280

# Row 140 (CUL)
$ws.Range("H140").Value = 13890620
$ws.Range("I140").Value = 20834414
$ws.Range("J140").Value = 3033.1667
$ws.Range("K140").Value = 62503242
$ws.Range("L140").Value = 9099.500100000001
$ws.Range("M140").Value = -62498062
$ws.Range("N140").Value = -19459.5001

$ws = $wb.Worksheets.Item("GSM")
# Row 132 (GSM)
$ws.Range("H132").Value = 6935.5
$ws.Range("I132").Value = 1650.8572
$ws.Range("K132").Value = 4952.571599999999
$ws.Range("M132").Value = -2422.571599999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 2500
$ws.Range("I7").Value = 2500
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2388
$ws.Range("N7").ClearContents()

# Row 46 (LTW)
$ws.Range("H46").Value = 2452535.5
$ws.Range("I46").Value = 3473168.5
$ws.Range("J46").Value = 3016.4
$ws.Range("K46").Value = 3473168.5
$ws.Range("L46").Value = 3016.4
$ws.Range("M46").Value = -3472980.5
$ws.Range("N46").Value = -3392.4

# Row 100 (LTW)
$ws.Range("H100").Value = 2279.75
$ws.Range("I100").Value = 2382.5715
$ws.Range("J100").Value = 2135.8
$ws.Range("K100").Value = 2382.5715
$ws.Range("L100").Value = 2135.8
$ws.Range("M100").Value = -1841.5715
$ws.Range("N100").Value = -3217.8

# Row 126 (LTW)
$ws.Range("H126").Value = 2500
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5030
$ws.Range("N126").ClearContents()

# Row 132 (LTW)
$ws.Range("H132").Value = 43966096
$ws.Range("I132").Value = 87913464
$ws.Range("K132").Value = 263740392
$ws.Range("M132").Value = -263737862

# Row 136 (LTW)
$ws.Range("H136").Value = 117916810
$ws.Range("I136").Value = 113556080
$ws.Range("J136").Value = 125002984
$ws.Range("K136").Value = 340668240
$ws.Range("L136").Value = 375008952
$ws.Range("M136").Value = -340665690
$ws.Range("N136").Value = -375014052

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (WVR)
$ws.Range("H132").Value = 28524.83
$ws.Range("I132").Value = 44892.176
$ws.Range("K132").Value = 134676.528
$ws.Range("M132").Value = -132146.528

# Row 136 (WVR)
$ws.Range("H136").Value = 9093052
$ws.Range("I136").Value = 20002452
$ws.Range("J136").Value = 1884.6666
$ws.Range("K136").Value = 60007356
$ws.Range("L136").Value = 5653.9998
$ws.Range("M136").Value = -60004806
$ws.Range("N136").Value = -10753.9998
